# PM04 Tidsregistrering for Toke.xlsx - add three new time-tracking entries
# (rows 22-24 on "Ark1") and move the active selection to E25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 22 - ATD01+02, 15/05/2020, 08:30 -> 09:15
$ws.Range("A22").Value = "ATD01+02"
$ws.Range("C22").Value = 43966
$ws.Range("D22").Value = 0.354166666666667
$ws.Range("E22").Value = 0.385416666666667

# Row 23 - OC0101+0102, 15/05/2020, 09:15 -> 13:00
$ws.Range("A23").Value = "OC0101+0102"
$ws.Range("C23").Value = 43966
$ws.Range("D23").Value = 0.385416666666667
$ws.Range("E23").Value = 0.541666666666667

# Row 24 - Rapport - Finansiering, 15/05/2020, 13:00 -> 16:00
$ws.Range("A24").Value = "Rapport – Finansiering"
$ws.Range("C24").Value = 43966
$ws.Range("D24").Value = 0.541666666666667
$ws.Range("E24").Value = 0.666666666666667

# Move the active selection on the sheet to E25 (matches the saved cursor
# position in the workbook after the edits).
$ws.Range("E25").Select()
